{"js": "// Replace the multiplication problems in the worksheet table with a new\n// set of problems. Each old expression is unique in the document, so a\n// straightforward search-and-replace keyed on the exact text is safe.\nconst replacements = [\n  [\"567\u00d78=\", \"391\u00d72=\"],\n  [\"882\u00d76=\", \"506\u00d78=\"],\n  [\"793\u00d72=\", \"676\u00d77=\"],\n  [\"643\u00d78=\", \"836\u00d79=\"],\n  [\"449\u00d72=\", \"846\u00d77=\"],\n  [\"146\u00d75=\", \"836\u00d75=\"],\n  [\"961\u00d79=\", \"671\u00d76=\"],\n  [\"367\u00d77=\", \"679\u00d75=\"],\n  [\"645\u00d75=\", \"728\u00d78=\"],\n  [\"504\u00d76=\", \"532\u00d79=\"],\n  [\"433\u00d78=\", \"670\u00d72=\"],\n  [\"355\u00d72=\", \"562\u00d74=\"],\n  [\"857\u00d74=\", \"240\u00d78=\"],\n  [\"295\u00d76=\", \"119\u00d74=\"],\n  [\"634\u00d73=\", \"443\u00d72=\"],\n  [\"725\u00d73=\", \"531\u00d75=\"],\n  [\"860\u00d75=\", \"909\u00d73=\"],\n  [\"689\u00d72=\", \"382\u00d75=\"],\n  [\"384\u00d74=\", \"782\u00d73=\"],\n  [\"499\u00d73=\", \"649\u00d79=\"],\n  [\"330\u00d73=\", \"491\u00d72=\"],\n  [\"400\u00d75=\", \"718\u00d72=\"],\n  [\"101\u00d76=\", \"344\u00d73=\"],\n  [\"657\u00d76=\", \"815\u00d74=\"],\n  [\"113\u00d72=\", \"142\u00d77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication problems in the worksheet table with a new set\n# of problems. Each old expression appears exactly once in the document, so\n# Find/Replace keyed on the exact text is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"567\u00d78=\", \"391\u00d72=\"),\n    @(\"882\u00d76=\", \"506\u00d78=\"),\n    @(\"793\u00d72=\", \"676\u00d77=\"),\n    @(\"643\u00d78=\", \"836\u00d79=\"),\n    @(\"449\u00d72=\", \"846\u00d77=\"),\n    @(\"146\u00d75=\", \"836\u00d75=\"),\n    @(\"961\u00d79=\", \"671\u00d76=\"),\n    @(\"367\u00d77=\", \"679\u00d75=\"),\n    @(\"645\u00d75=\", \"728\u00d78=\"),\n    @(\"504\u00d76=\", \"532\u00d79=\"),\n    @(\"433\u00d78=\", \"670\u00d72=\"),\n    @(\"355\u00d72=\", \"562\u00d74=\"),\n    @(\"857\u00d74=\", \"240\u00d78=\"),\n    @(\"295\u00d76=\", \"119\u00d74=\"),\n    @(\"634\u00d73=\", \"443\u00d72=\"),\n    @(\"725\u00d73=\", \"531\u00d75=\"),\n    @(\"860\u00d75=\", \"909\u00d73=\"),\n    @(\"689\u00d72=\", \"382\u00d75=\"),\n    @(\"384\u00d74=\", \"782\u00d73=\"),\n    @(\"499\u00d73=\", \"649\u00d79=\"),\n    @(\"330\u00d73=\", \"491\u00d72=\"),\n    @(\"400\u00d75=\", \"718\u00d72=\"),\n    @(\"101\u00d76=\", \"344\u00d73=\"),\n    @(\"657\u00d76=\", \"815\u00d74=\"),\n    @(\"113\u00d72=\", \"142\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2)\n}\n"}
